$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename rain model constants
$ws.Range("A19").Value = "fi_lidar_rain_reflectivity"
$ws.Range("A20").Value = "fi_lidar_rain_intensity"

# Move selection to A20 (last edited cell)
$ws.Range("A20").Select()
